# Updates crypto price/volume data per the Wed Nov  6 05:41:36 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.434.53"
$ws.Range("E2").Value = "  +8.65%  "
$ws.Range("D3").Value = "2.587.34"
$ws.Range("E3").Value = "  +6.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.84"
$ws.Range("E5").Value = "  +15.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "579.83"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +4.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.205"
$ws.Range("E9").Value = "  +25.09%  "
$ws.Range("D10").Value = "2.585.09"
$ws.Range("E10").Value = "  +6.48%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("E12").Value = "  +8.52%  "
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000191"
$ws.Range("E14").Value = "  +9.85%  "
$ws.Range("D15").Value = "74.210.47"
$ws.Range("E15").Value = "  +8.51%  "
$ws.Range("D16").Value = "3.046.89"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.16"
$ws.Range("E17").Value = "  +13.32%  "
$ws.Range("D18").Value = "2.593.11"
$ws.Range("E18").Value = "  +6.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.95"
$ws.Range("E19").Value = "  +30.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  +12.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.51"
$ws.Range("E21").Value = "  +12.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("E22").Value = "  +19.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.05"
$ws.Range("E23").Value = "  +6.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.65"
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.14"
$ws.Range("E26").Value = "  +12.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("E27").Value = "  +11.84%  "
$ws.Range("D28").Value = "2.716.61"
$ws.Range("E28").Value = "  +6.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "0.0₃0941"
$ws.Range("E30").Value = "  +14.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.95"
$ws.Range("E31").Value = "  +11.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "503.64"
$ws.Range("E32").Value = "  +18.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  +18.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.71"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  +12.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.83"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.19"
$ws.Range("E38").Value = "  +7.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.39"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.91"
$ws.Range("E41").Value = "  +13.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.68"
$ws.Range("E42").Value = "  +12.51%  "
$ws.Range("E43").Value = "  +7.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.44"
$ws.Range("E44").Value = "  +19.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.02"
$ws.Range("E45").Value = "  +4.51%  "
$ws.Range("E46").Value = "  +7.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.70"
$ws.Range("E47").Value = "  +13.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0816"
$ws.Range("E48").Value = "  +14.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.61"
$ws.Range("E49").Value = "  +8.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.521"
$ws.Range("E50").Value = "  +8.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0969"
$ws.Range("E51").Value = "  +5.95%  "
